$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.539.11"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.580.96"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'582.38"
$ws.Range("D6").Value = "'165.89"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "2.580.69"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "'26.79"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "3.051.23"
$ws.Range("D17").Value = "66.377.87"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "2.584.17"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "'11.43"
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").Value = "'7.75"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "'352.17"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("D26").Value = "'10.04"
$ws.Range("E26").Value = "  -7.93%  "
$ws.Range("D27").Value = "'69.02"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "2.711.55"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.0₃0989"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "'536.44"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").Value = "'8.01"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").Value = "'156.78"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").Value = "'18.76"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'5.12"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'2.42"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").Value = "'149.18"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'1.70"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "'0.0760"
$ws.Range("E51").Value = "  -1.49%  "
